# Natmi following Dr Hou advice
# Update the Adam15-Itgb3 LR-pairs sheet: every sending/target cluster
# combination now aggregates 3 samples (E/K column counts) instead of 1,
# and the dependent numeric columns (G..T) are recomputed accordingly.
# Also adds the 5 missing rows for the sCs sending cluster (rows 22-26),
# completing the full 5x5 ECs/FAPs/M1/M2/sCs sending x target grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A="ECs"; B="Adam15"; C="Itgb3"; D="ECs"; E=3; F=1; G=69.00849933333333; H=207.025498; I=0.4252801448282221; J=0.425280144828222; K=3; L=1; M=6.610178333333334; N=19.830535; O=0.6561418030098394; P=0.6561418030098395; Q=456.1584871090478; R=4105.42638398143; S=0.2790440810118752; T=0.2790440810118753 },
    @{ A="ECs"; B="Adam15"; C="Itgb3"; D="FAPs"; E=3; F=1; G=69.00849933333333; H=207.025498; I=0.4252801448282221; J=0.425280144828222; K=3; L=1; M=2.754304; N=8.262912; O=0.2733986742057961; P=0.2733986742057962; Q=190.0703857477973; R=1710.633471730176; S=0.1162710277620849; T=0.1162710277620849 },
    @{ A="ECs"; B="Adam15"; C="Itgb3"; D="M1"; E=3; F=1; G=69.00849933333333; H=207.025498; I=0.4252801448282221; J=0.425280144828222; K=3; L=1; M=0.2292836666666667; N=0.687851; O=0.02275923445041301; P=0.02275923445041302; Q=15.82252175831089; R=142.402695824798; S=0.009679050523251107; T=0.009679050523251107 },
    @{ A="ECs"; B="Adam15"; C="Itgb3"; D="M2"; E=3; F=1; G=69.00849933333333; H=207.025498; I=0.4252801448282221; J=0.425280144828222; K=2; L=0.6666666666666666; M=0.214866; N=0.644598; O=0.02132810304596101; P=0.02132810304596101; Q=14.827580217756; R=133.448221959804; S=0.009070418752297542; T=0.009070418752297544 },
    @{ A="ECs"; B="Adam15"; C="Itgb3"; D="sCs"; E=3; F=1; G=69.00849933333333; H=207.025498; I=0.4252801448282221; J=0.425280144828222; K=3; L=1; M=0.2656816666666666; N=0.797045; O=0.02637218528799033; P=0.02637218528799034; Q=18.33429311704555; R=165.00863805341; S=0.01121556677871324; T=0.01121556677871324 },
    @{ A="FAPs"; B="Adam15"; C="Itgb3"; D="ECs"; E=3; F=1; G=14.694925; H=44.084775; I=0.09056072647012584; J=0.09056072647012582; K=3; L=1; M=6.610178333333334; N=19.830535; O=0.6561418030098394; P=0.6561418030098395; Q=97.13607484495833; R=874.2246736046251; S=0.05942067834798925; T=0.05942067834798925 },
    @{ A="FAPs"; B="Adam15"; C="Itgb3"; D="FAPs"; E=3; F=1; G=14.694925; H=44.084775; I=0.09056072647012584; J=0.09056072647012582; K=3; L=1; M=2.754304; N=8.262912; O=0.2733986742057961; P=0.2733986742057962; Q=40.4742907072; R=364.2686163648; S=0.02475918255204615; T=0.02475918255204616 },
    @{ A="FAPs"; B="Adam15"; C="Itgb3"; D="M1"; E=3; F=1; G=14.694925; H=44.084775; I=0.09056072647012584; J=0.09056072647012582; K=3; L=1; M=0.2292836666666667; N=0.687851; O=0.02275923445041301; P=0.02275923445041302; Q=3.369306285391667; R=30.323756568525; S=0.002061092805733317; T=0.002061092805733317 },
    @{ A="FAPs"; B="Adam15"; C="Itgb3"; D="M2"; E=3; F=1; G=14.694925; H=44.084775; I=0.09056072647012584; J=0.09056072647012582; K=2; L=0.6666666666666666; M=0.214866; N=0.644598; O=0.02132810304596101; P=0.02132810304596101; Q=3.15743975505; R=28.41695779545; S=0.001931488506071933; T=0.001931488506071933 },
    @{ A="FAPs"; B="Adam15"; C="Itgb3"; D="sCs"; E=3; F=1; G=14.694925; H=44.084775; I=0.09056072647012584; J=0.09056072647012582; K=3; L=1; M=0.2656816666666666; N=0.797045; O=0.02637218528799033; P=0.02637218528799034; Q=3.904172165541666; R=35.137549489875; S=0.002388284258285169; T=0.002388284258285169 },
    @{ A="M1"; B="Adam15"; C="Itgb3"; D="ECs"; E=3; F=1; G=43.213838; H=129.641514; I=0.2663148374586689; J=0.2663148374586689; K=3; L=1; M=6.610178333333334; N=19.830535; O=0.6561418030098394; P=0.6561418030098395; Q=285.6511756477767; R=2570.86058082999; S=0.1747402976184033; T=0.1747402976184033 },
    @{ A="M1"; B="Adam15"; C="Itgb3"; D="FAPs"; E=3; F=1; G=43.213838; H=129.641514; I=0.2663148374586689; J=0.2663148374586689; K=3; L=1; M=2.754304; N=8.262912; O=0.2733986742057961; P=0.2733986742057962; Q=119.024046858752; R=1071.216421728768; S=0.07281012348253217; T=0.07281012348253219 },
    @{ A="M1"; B="Adam15"; C="Itgb3"; D="M1"; E=3; F=1; G=43.213838; H=129.641514; I=0.2663148374586689; J=0.2663148374586689; K=3; L=1; M=0.2292836666666667; N=0.687851; O=0.02275923445041301; P=0.02275923445041302; Q=9.908227227379333; R=89.17404504641399; S=0.00606112182334548; T=0.006061121823345479 },
    @{ A="M1"; B="Adam15"; C="Itgb3"; D="M2"; E=3; F=1; G=43.213838; H=129.641514; I=0.2663148374586689; J=0.2663148374586689; K=2; L=0.6666666666666666; M=0.214866; N=0.644598; O=0.02132810304596101; P=0.02132810304596101; Q=9.285184515708; R=83.566660641372; S=0.005679990295986848; T=0.005679990295986847 },
    @{ A="M1"; B="Adam15"; C="Itgb3"; D="sCs"; E=3; F=1; G=43.213838; H=129.641514; I=0.2663148374586689; J=0.2663148374586689; K=3; L=1; M=0.2656816666666666; N=0.797045; O=0.02637218528799033; P=0.02637218528799034; Q=11.48112450290333; R=103.33012052613; S=0.007023304238401045; T=0.007023304238401046 },
    @{ A="M2"; B="Adam15"; C="Itgb3"; D="ECs"; E=3; F=1; G=34.20250066666667; H=102.607502; I=0.2107804774415859; J=0.2107804774415858; K=3; L=1; M=6.610178333333334; N=19.830535; O=0.6561418030098394; P=0.6561418030098395; Q=226.0846288526189; R=2034.76165967357; S=0.138301882507797; T=0.1383018825077969 },
    @{ A="M2"; B="Adam15"; C="Itgb3"; D="FAPs"; E=3; F=1; G=34.20250066666667; H=102.607502; I=0.2107804774415859; J=0.2107804774415858; K=3; L=1; M=2.754304; N=8.262912; O=0.2733986742057961; P=0.2733986742057962; Q=94.20408439620266; R=847.836759565824; S=0.0576271030809943; T=0.05762710308099431 },
    @{ A="M2"; B="Adam15"; C="Itgb3"; D="M1"; E=3; F=1; G=34.20250066666667; H=102.607502; I=0.2107804774415859; J=0.2107804774415858; K=3; L=1; M=0.2292836666666667; N=0.687851; O=0.02275923445041301; P=0.02275923445041302; Q=7.842074762022444; R=70.578672858202; S=0.004797202303663045; T=0.004797202303663044 },
    @{ A="M2"; B="Adam15"; C="Itgb3"; D="M2"; E=3; F=1; G=34.20250066666667; H=102.607502; I=0.2107804774415859; J=0.2107804774415858; K=2; L=0.6666666666666666; M=0.214866; N=0.644598; O=0.02132810304596101; P=0.02132810304596101; Q=7.348954508244; R=66.140590574196; S=0.004495547742951004; T=0.004495547742951004 },
    @{ A="M2"; B="Adam15"; C="Itgb3"; D="sCs"; E=3; F=1; G=34.20250066666667; H=102.607502; I=0.2107804774415859; J=0.2107804774415858; K=3; L=1; M=0.2656816666666666; N=0.797045; O=0.02637218528799033; P=0.02637218528799034; Q=9.086977381287777; R=81.78279643159; S=0.00555874180618057; T=0.00555874180618057 },
    @{ A="sCs"; B="Adam15"; C="Itgb3"; D="ECs"; E=3; F=1; G=1.146216666666667; H=3.43865; I=0.007063813801397381; J=0.007063813801397379; K=3; L=1; M=6.610178333333334; N=19.830535; O=0.6561418030098394; P=0.6561418030098395; Q=7.576696575305556; R=68.19026917775; S=0.004634863523774665; T=0.004634863523774665 },
    @{ A="sCs"; B="Adam15"; C="Itgb3"; D="FAPs"; E=3; F=1; G=1.146216666666667; H=3.43865; I=0.007063813801397381; J=0.007063813801397379; K=3; L=1; M=2.754304; N=8.262912; O=0.2733986742057961; P=0.2733986742057962; Q=3.157029149866667; R=28.4132623488; S=0.001931237328138649; T=0.001931237328138649 },
    @{ A="sCs"; B="Adam15"; C="Itgb3"; D="M1"; E=3; F=1; G=1.146216666666667; H=3.43865; I=0.007063813801397381; J=0.007063813801397379; K=3; L=1; M=0.2292836666666667; N=0.687851; O=0.02275923445041301; P=0.02275923445041302; Q=0.2628087601277778; R=2.36527884115; S=0.0001607669944200662; T=0.0001607669944200662 },
    @{ A="sCs"; B="Adam15"; C="Itgb3"; D="M2"; E=3; F=1; G=1.146216666666667; H=3.43865; I=0.007063813801397381; J=0.007063813801397379; K=2; L=0.6666666666666666; M=0.214866; N=0.644598; O=0.02132810304596101; P=0.02132810304596101; Q=0.2462829903; R=2.2165469127; S=0.0001506577486536849; T=0.0001506577486536849 },
    @{ A="sCs"; B="Adam15"; C="Itgb3"; D="sCs"; E=3; F=1; G=1.146216666666667; H=3.43865; I=0.007063813801397381; J=0.007063813801397379; K=3; L=1; M=0.2656816666666666; N=0.797045; O=0.02637218528799033; P=0.02637218528799034; Q=0.3045287543611111; R=2.74075878925; S=0.0001862882064103151; T=0.0001862882064103151 }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $r++
}
